$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CIPLItem")

# --- ReferenceNo column (A2:A4): 4100184111 -> 4100158785 ---
# Setting .Value directly on a numeric cell resets any existing cell style
# (the original cells carry a quote-prefix style), so re-stamp the format
# from a still-untouched sibling cell right after writing the new value.
$ws.Range("A2").Value = 4100158785
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A3").Value = 4100158785
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A4").Value = 4100158785
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# --- GrossWeigth column (S2:S4): 23.000 -> 54.001 (stays text) ---
$ws.Range("S2").NumberFormat = "@"
$ws.Range("S2").Value = "54.001"
$ws.Range("D2").Copy() | Out-Null
$ws.Range("S2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("S3").NumberFormat = "@"
$ws.Range("S3").Value = "54.001"
$ws.Range("D2").Copy() | Out-Null
$ws.Range("S3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("S4").NumberFormat = "@"
$ws.Range("S4").Value = "54.001"
$ws.Range("D2").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# --- View/selection state ---
$ws.Range("B4").Select()
